# Apply "Final Update, Documentation and Mastering" changes to Asset_List.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update Status column (D9:D12) from "X" to "Done" ---
$ws.Range("D9").Value = "Done"
$ws.Range("D10").Value = "Done"
$ws.Range("D11").Value = "Done"
$ws.Range("D12").Value = "Done"

# --- Add new row 15: Title_Ambience (added before the Background description edit so new
#     shared strings land in the same order as the authored workbook) ---
$ws.Range("A15").Value = "Title_Ambience"
$ws.Range("B15").Value = "Ambience"
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = "Done"
$ws.Range("E15").Value = "Scatterer instrument that plays in the title screen"

# --- Update Background description: "Louder" -> "Quieter" ---
$ws.Range("E12").Value = "Background ambience that plays throughout the game. Quieter during menus, contains bird chirps, wind whistling, dog bark, windchimes, bugs buzzing, etc"

# --- Adjust column widths: column A gets its own (wider) width, column B keeps the old shared width ---
$ws.Columns.Item(1).ColumnWidth = 13.6640625
$ws.Columns.Item(2).ColumnWidth = 11.33203125

# --- Update sheet view: scroll + selection ---
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("E12").Select()
